$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.102.57'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').Value = '2.897.07'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '''590.36'
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').Value = '''141.34'
$ws.Range('E6').Value = '  -3.19%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '''0.494'
$ws.Range('E8').Value = '  -2.48%  '
$ws.Range('D9').Value = '''6.90'
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('D10').Value = '''0.138'
$ws.Range('E10').Value = '  -3.93%  '
$ws.Range('D11').Value = '''0.430'
$ws.Range('E11').Value = '  -2.28%  '
$ws.Range('D12').Value = '''0.0000218'
$ws.Range('E12').Value = '  -3.39%  '
$ws.Range('D13').Value = '''32.45'
$ws.Range('E13').Value = '  -3.25%  '
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').Value = '3.390.57'
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('D16').Value = '61.114.52'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').Value = '2.908.14'
$ws.Range('E17').Value = '  -0.20%  '
$ws.Range('D18').Value = '''6.52'
$ws.Range('E18').Value = '  -2.46%  '
$ws.Range('D19').Value = '''427.97'
$ws.Range('E19').Value = '  -0.68%  '
$ws.Range('D20').Value = '''13.37'
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('D21').Value = '''0.656'
$ws.Range('E21').Value = '  -3.12%  '
$ws.Range('D22').Value = '''6.94'
$ws.Range('E22').Value = '  -1.14%  '
$ws.Range('D23').Value = '''80.74'
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('D24').Value = '''10.59'
$ws.Range('E24').Value = '  -3.80%  '
$ws.Range('D25').Value = '''2.08'
$ws.Range('E25').Value = '  -4.89%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').Value = '''11.47'
$ws.Range('E27').Value = '  -2.94%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').Value = '''2.55'
$ws.Range('E28').Value = '  -2.23%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').Value = '''2.11'
$ws.Range('E29').Value = '  -7.09%  '
$ws.Range('D30').Value = '''6.77'
$ws.Range('E30').Value = '  -3.57%  '
$ws.Range('D31').Value = '''26.19'
$ws.Range('E31').Value = '  -1.46%  '
$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D32').Value = '''1.00'
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '''0.106'
$ws.Range('E33').Value = '  -2.12%  '
$ws.Range('D34').Value = '0.0₃0858'
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('D35').Value = '''0.980'
$ws.Range('E35').Value = '  -2.79%  '
$ws.Range('D36').Value = '''5.46'
$ws.Range('E36').Value = '  -2.70%  '
$ws.Range('D37').Value = '''2.85'
$ws.Range('E37').Value = '  -5.38%  '
$ws.Range('D38').Value = '''1.93'
$ws.Range('E38').Value = '  -2.29%  '
$ws.Range('D39').Value = '''0.119'
$ws.Range('E39').Value = '  -2.18%  '
$ws.Range('D40').Value = '''8.31'
$ws.Range('E40').Value = '  -2.76%  '
$ws.Range('D41').Value = '''40.80'
$ws.Range('E41').Value = '  +0.34%  '
$ws.Range('D42').Value = '''0.268'
$ws.Range('E42').Value = '  -6.07%  '
$ws.Range('D43').Value = '2.660.02'
$ws.Range('E43').Value = '  -1.21%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value = '''132.86'
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '''0.0335'
$ws.Range('E45').Value = '  -2.25%  '
$ws.Range('D46').Value = '''347.78'
$ws.Range('E46').Value = '  -8.15%  '
$ws.Range('D48').Value = '''22.82'
$ws.Range('E48').Value = '  -4.33%  '
$ws.Range('D49').Value = '''0.103'
$ws.Range('E49').Value = '  -2.80%  '
$ws.Range('D50').Value = '''1.95'
$ws.Range('E50').Value = '  -2.26%  '
$ws.Range('D51').Value = '''0.123'
$ws.Range('E51').Value = '  -0.82%  '
